$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.4047206666666667
$ws.Range("H2").Value = 1.214162
$ws.Range("I2").Value = 0.02578034976888792
$ws.Range("J2").Value = 0.02578034976888792
$ws.Range("M2").Value = 133.7780026666667
$ws.Range("N2").Value = 401.334008
$ws.Range("O2").Value = 0.50863533211804
$ws.Range("P2").Value = 0.5086353321180399
$ws.Range("Q2").Value = 54.14272242458845
$ws.Range("R2").Value = 487.284501821296
$ws.Range("S2").Value = 0.01311279676681754
$ws.Range("T2").Value = 0.01311279676681754

$ws.Range("G3").Value = 0.4047206666666667
$ws.Range("H3").Value = 1.214162
$ws.Range("I3").Value = 0.02578034976888792
$ws.Range("J3").Value = 0.02578034976888792
$ws.Range("O3").Value = 0.1993888292903622
$ws.Range("P3").Value = 0.1993888292903622
$ws.Range("Q3").Value = 21.224349464434
$ws.Range("R3").Value = 191.019145179906
$ws.Range("S3").Value = 0.005140313759114622
$ws.Range("T3").Value = 0.005140313759114622

$ws.Range("G4").Value = 0.4047206666666667
$ws.Range("H4").Value = 1.214162
$ws.Range("I4").Value = 0.02578034976888792
$ws.Range("J4").Value = 0.02578034976888792
$ws.Range("M4").Value = 21.197691
$ws.Range("N4").Value = 63.593073
$ws.Range("O4").Value = 0.08059542216956049
$ws.Range("P4").Value = 0.08059542216956046
$ws.Range("Q4").Value = 8.579143633314001
$ws.Range("R4").Value = 77.21229269982601
$ws.Range("S4").Value = 0.002077778173302453
$ws.Range("T4").Value = 0.002077778173302452

$ws.Range("G5").Value = 0.4047206666666667
$ws.Range("H5").Value = 1.214162
$ws.Range("I5").Value = 0.02578034976888792
$ws.Range("J5").Value = 0.02578034976888792
$ws.Range("M5").Value = 55.59592133333333
$ws.Range("N5").Value = 166.787764
$ws.Range("O5").Value = 0.2113804164220374
$ws.Range("P5").Value = 0.2113804164220373
$ws.Range("Q5").Value = 22.50081834597422
$ws.Range("R5").Value = 202.507365113768
$ws.Range("S5").Value = 0.005449461069653303
$ws.Range("T5").Value = 0.005449461069653301

$ws.Range("I6").Value = 0.03222381288358415
$ws.Range("J6").Value = 0.03222381288358415
$ws.Range("M6").Value = 133.7780026666667
$ws.Range("N6").Value = 401.334008
$ws.Range("O6").Value = 0.50863533211804
$ws.Range("P6").Value = 0.5086353321180399
$ws.Range("Q6").Value = 67.67499169166756
$ws.Range("R6").Value = 609.0749252250081
$ws.Range("S6").Value = 0.0163901697681514
$ws.Range("T6").Value = 0.0163901697681514

$ws.Range("I7").Value = 0.03222381288358415
$ws.Range("J7").Value = 0.03222381288358415
$ws.Range("O7").Value = 0.1993888292903622
$ws.Range("P7").Value = 0.1993888292903622
$ws.Range("S7").Value = 0.006425068326129535
$ws.Range("T7").Value = 0.006425068326129535

$ws.Range("I8").Value = 0.03222381288358415
$ws.Range("J8").Value = 0.03222381288358415
$ws.Range("M8").Value = 21.197691
$ws.Range("N8").Value = 63.593073
$ws.Range("O8").Value = 0.08059542216956049
$ws.Range("P8").Value = 0.08059542216956046
$ws.Range("Q8").Value = 10.723389000522
$ws.Range("R8").Value = 96.51050100469801
$ws.Range("S8").Value = 0.002597091803265387
$ws.Range("T8").Value = 0.002597091803265387

$ws.Range("I9").Value = 0.03222381288358415
$ws.Range("J9").Value = 0.03222381288358415
$ws.Range("M9").Value = 55.59592133333333
$ws.Range("N9").Value = 166.787764
$ws.Range("O9").Value = 0.2113804164220374
$ws.Range("P9").Value = 0.2113804164220373
$ws.Range("Q9").Value = 28.12460523647378
$ws.Range("R9").Value = 253.121447128264
$ws.Range("S9").Value = 0.00681148298603783
$ws.Range("T9").Value = 0.006811482986037829

$ws.Range("G10").Value = 1.039987
$ws.Range("H10").Value = 3.119961
$ws.Range("I10").Value = 0.06624625531460326
$ws.Range("J10").Value = 0.06624625531460326
$ws.Range("M10").Value = 133.7780026666667
$ws.Range("N10").Value = 401.334008
$ws.Range("O10").Value = 0.50863533211804
$ws.Range("P10").Value = 0.5086353321180399
$ws.Range("Q10").Value = 139.1273836592987
$ws.Range("R10").Value = 1252.146452933688
$ws.Range("S10").Value = 0.03369518607351971
$ws.Range("T10").Value = 0.0336951860735197

$ws.Range("G11").Value = 1.039987
$ws.Range("H11").Value = 3.119961
$ws.Range("I11").Value = 0.06624625531460326
$ws.Range("J11").Value = 0.06624625531460326
$ws.Range("O11").Value = 0.1993888292903622
$ws.Range("P11").Value = 0.1993888292903622
$ws.Range("Q11").Value = 54.538968094377
$ws.Range("R11").Value = 490.8507128493931
$ws.Range("S11").Value = 0.01320876329204918
$ws.Range("T11").Value = 0.01320876329204918

$ws.Range("G12").Value = 1.039987
$ws.Range("H12").Value = 3.119961
$ws.Range("I12").Value = 0.06624625531460326
$ws.Range("J12").Value = 0.06624625531460326
$ws.Range("M12").Value = 21.197691
$ws.Range("N12").Value = 63.593073
$ws.Range("O12").Value = 0.08059542216956049
$ws.Range("P12").Value = 0.08059542216956046
$ws.Range("Q12").Value = 22.045323070017
$ws.Range("R12").Value = 198.407907630153
$ws.Range("S12").Value = 0.005339144914232941
$ws.Range("T12").Value = 0.005339144914232938

$ws.Range("G13").Value = 1.039987
$ws.Range("H13").Value = 3.119961
$ws.Range("I13").Value = 0.06624625531460326
$ws.Range("J13").Value = 0.06624625531460326
$ws.Range("M13").Value = 55.59592133333333
$ws.Range("N13").Value = 166.787764
$ws.Range("O13").Value = 0.2113804164220374
$ws.Range("P13").Value = 0.2113804164220373
$ws.Range("Q13").Value = 57.81903543968933
$ws.Range("R13").Value = 520.3713189572039
$ws.Range("S13").Value = 0.01400316103480144
$ws.Range("T13").Value = 0.01400316103480144

$ws.Range("G14").Value = 13.74822133333333
$ws.Range("H14").Value = 41.244664
$ws.Range("I14").Value = 0.8757495820329246
$ws.Range("J14").Value = 0.8757495820329247
$ws.Range("M14").Value = 133.7780026666667
$ws.Range("N14").Value = 401.334008
$ws.Range("O14").Value = 0.50863533211804
$ws.Range("P14").Value = 0.5086353321180399
$ws.Range("Q14").Value = 1839.20959019259
$ws.Range("R14").Value = 16552.88631173331
$ws.Range("S14").Value = 0.4454371795095514
$ws.Range("T14").Value = 0.4454371795095513

$ws.Range("G15").Value = 13.74822133333333
$ws.Range("H15").Value = 41.244664
$ws.Range("I15").Value = 0.8757495820329246
$ws.Range("J15").Value = 0.8757495820329247
$ws.Range("O15").Value = 0.1993888292903622
$ws.Range("P15").Value = 0.1993888292903622
$ws.Range("Q15").Value = 720.9838244642481
$ws.Range("R15").Value = 6488.854420178232
$ws.Range("S15").Value = 0.1746146839130689
$ws.Range("T15").Value = 0.1746146839130689

$ws.Range("G16").Value = 13.74822133333333
$ws.Range("H16").Value = 41.244664
$ws.Range("I16").Value = 0.8757495820329246
$ws.Range("J16").Value = 0.8757495820329247
$ws.Range("M16").Value = 21.197691
$ws.Range("N16").Value = 63.593073
$ws.Range("O16").Value = 0.08059542216956049
$ws.Range("P16").Value = 0.08059542216956046
$ws.Range("Q16").Value = 291.430547623608
$ws.Range("R16").Value = 2622.874928612472
$ws.Range("S16").Value = 0.07058140727875971
$ws.Range("T16").Value = 0.07058140727875969

$ws.Range("G17").Value = 13.74822133333333
$ws.Range("H17").Value = 41.244664
$ws.Range("I17").Value = 0.8757495820329246
$ws.Range("J17").Value = 0.8757495820329247
$ws.Range("M17").Value = 55.59592133333333
$ws.Range("N17").Value = 166.787764
$ws.Range("O17").Value = 0.2113804164220374
$ws.Range("P17").Value = 0.2113804164220373
$ws.Range("Q17").Value = 764.3450317212551
$ws.Range("R17").Value = 6879.105285491295
$ws.Range("S17").Value = 0.1851163113315448
$ws.Range("T17").Value = 0.1851163113315447
